$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply new header style (bold/border) to newly introduced columns L1:O1 ---
$ws.Range("A1").Copy($ws.Range("L1:O1"))

# --- Row 1 ---
$ws.Range("A1").Value = "Best Estimator"
$ws.Range("B1").Value = "Best Score"
$ws.Range("C1").Value = "Best Params"
$ws.Range("D1").Value = "CV Train F1"
$ws.Range("E1").Value = "CV Test F1"
$ws.Range("F1").Value = "Validation F1"
$ws.Range("G1").Value = "CV Train Precision"
$ws.Range("H1").Value = "CV Test Precision"
$ws.Range("I1").Value = "Validation Precision"
$ws.Range("J1").Value = "CV Train Recall"
$ws.Range("K1").Value = "CV Test Recall"
$ws.Range("L1").Value = "Validation Recall"
$ws.Range("M1").Value = "Y Val (Validation)"
$ws.Range("N1").Value = "Y Pred (Validation)"
$ws.Range("O1").Value = "Seed"

# --- Row 2 ---
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fd587256220>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=1, class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=5, random_state=42))])"
$ws.Range("B2").Value = 0.73
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd587256250>, 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 1}"
$ws.Range("D2").Value = 0.7129496054804985
$ws.Range("E2").Value = 0.5679567876567877
$ws.Range("F2").Value = 0.7234042553191491
$ws.Range("G2").Value = 0.6479264075900082
$ws.Range("H2").Value = 0.5156095238095237
$ws.Range("I2").Value = 0.68
$ws.Range("J2").Value = 0.8366666666666667
$ws.Range("K2").Value = 0.6968
$ws.Range("L2").Value = 0.7727272727272727
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1 0 1 1 1 1 0 0 0 0 1 0 1]"
$ws.Range("N2").Value = "[1 1 1 0 1 1 1 1 1 1 1 1 1 0 1 0 0 0 0 1 1 1 1 1 0 1 1 1 0 1 0 0 0 1 1 1]"
$ws.Range("O2").Value = 42

# --- Row 3 ---
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fd5872569a0>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=5, class_weight='balanced',
                                                 random_state=42),
                                   random_state=42))])"
$ws.Range("B3").Value = 0.6771428571428572
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd587433b20>, 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 5}"
$ws.Range("D3").Value = 0.6871184070925034
$ws.Range("E3").Value = 0.5618914030414032
$ws.Range("F3").Value = 0.6363636363636365
$ws.Range("G3").Value = 0.6356049308725974
$ws.Range("H3").Value = 0.4969015873015873
$ws.Range("I3").Value = 0.7
$ws.Range("J3").Value = 0.8038999999999998
$ws.Range("K3").Value = 0.694
$ws.Range("L3").Value = 0.5833333333333334
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1]"
$ws.Range("N3").Value = "[0 1 1 1 1 0 0 0 1 1 0 0 1 0 1 0 1 1 1 1 1 0 1 1 1 1 0 0 0 0 0 0 1 1 0 1]"
$ws.Range("O3").Value = 69

# --- Row 4 ---
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fd587256be0>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=1, class_weight='balanced',
                                                 kernel='poly',
                                                 random_state=42),
                                   n_estimators=5, random_state=42))])"
$ws.Range("B4").Value = 0.6249999999999999
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd58b90b280>, 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__kernel': 'poly', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 1}"
$ws.Range("D4").Value = 0.5174177384007134
$ws.Range("E4").Value = 0.3690166722166722
$ws.Range("F4").Value = 0.4444444444444444
$ws.Range("G4").Value = 0.5279069628650789
$ws.Range("H4").Value = 0.357145238095238
$ws.Range("I4").Value = 0.8
$ws.Range("J4").Value = 0.5622105263157895
$ws.Range("K4").Value = 0.4223999999999999
$ws.Range("L4").Value = 0.3076923076923077
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 0 1 0 1 1 1 0 1]"
$ws.Range("N4").Value = "[0 1 1 0 0 1 0 1 1 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 1 0 0 0 1 0 0 1 0 1 0]"
$ws.Range("O4").Value = 23

# --- Row 5 ---
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fd587256490>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=1, class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=5, random_state=42))])"
$ws.Range("B5").Value = 0.6452380952380953
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd58b90b670>, 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 1}"
$ws.Range("D5").Value = 0.6758508792164776
$ws.Range("E5").Value = 0.5351280497280498
$ws.Range("F5").Value = 0.6521739130434783
$ws.Range("G5").Value = 0.6146296692144559
$ws.Range("H5").Value = 0.466306746031746
$ws.Range("I5").Value = 0.625
$ws.Range("J5").Value = 0.796452380952381
$ws.Range("K5").Value = 0.6739999999999999
$ws.Range("L5").Value = 0.6818181818181818
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 1 1 1 1 1 0]"
$ws.Range("N5").Value = "[0 1 0 0 1 0 1 1 0 1 1 1 0 0 1 0 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1]"
$ws.Range("O5").Value = 99

# --- Row 6 ---
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7fd58b90b220>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=5, class_weight='balanced',
                                                 kernel='linear',
                                                 random_state=42),
                                   random_state=42))])"
$ws.Range("B6").Value = 0.75
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd58723a9d0>, 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__kernel': 'linear', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 5}"
$ws.Range("D6").Value = 0.7448237877459697
$ws.Range("E6").Value = 0.5918844322344323
$ws.Range("F6").Value = 0.6511627906976744
$ws.Range("G6").Value = 0.6800642468822288
$ws.Range("H6").Value = 0.5056718253968254
$ws.Range("I6").Value = 0.6086956521739131
$ws.Range("J6").Value = 0.8579545454545454
$ws.Range("K6").Value = 0.7509999999999999
$ws.Range("L6").Value = 0.7
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1 1 0 1 0 1 1 1 1 1 1 1 0]"
$ws.Range("N6").Value = "[1 1 1 0 1 1 1 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1 1 1 0 0 1 1 0 1 1 1 1 0 0 0]"
$ws.Range("O6").Value = 89

Write-Host "Edit applied successfully"